$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row = 14

$ws.Cells.Item($row, 1).Value = "'2023-06-01"
$ws.Cells.Item($row, 2).Value = "'19:43:03"
$ws.Cells.Item($row, 3).Value = "Thursday"
$ws.Cells.Item($row, 4).Value = "'22"

$ws.Cells.Item($row, 1).ClearFormats()
$ws.Cells.Item($row, 2).ClearFormats()
$ws.Cells.Item($row, 3).ClearFormats()
$ws.Cells.Item($row, 4).ClearFormats()

$ws.Cells.Item($row, 5).Value = 120516
$ws.Cells.Item($row, 6).Value = 133763
$ws.Cells.Item($row, 7).Value = 158220
$ws.Cells.Item($row, 8).Value = 129338
$ws.Cells.Item($row, 9).Value = 173929
$ws.Cells.Item($row, 10).Value = 111686
$ws.Cells.Item($row, 11).Value = 198924
$ws.Cells.Item($row, 12).Value = 217494
$ws.Cells.Item($row, 13).Value = 170584
$ws.Cells.Item($row, 14).Value = 118211
$ws.Cells.Item($row, 15).Value = 37859
$ws.Cells.Item($row, 16).Value = 34930
$ws.Cells.Item($row, 17).Value = 49946
$ws.Cells.Item($row, 18).Value = -1
$ws.Cells.Item($row, 19).Value = 36733
$ws.Cells.Item($row, 20).Value = -1
